# removed seconds handling in search widgets (ngeo 366)
# Row 8 ("seconds handling in dates widgets") is removed; the following row
# (map display / In 2D mode...) moves up to take its place (row 8), and the
# sheet shrinks by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Delete the entire row 8 (the "seconds handling" issue). This shifts row 9
# (and its formatting, including row height) up into row 8's position.
$ws.Rows.Item(8).Delete()

# Update the active selection to match the target state.
$ws.Range("C2").Select()
